$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Value = "2025/12/04 00:00"
$ws.Range("B32").Value = "-"
$ws.Range("C32").Value = "-"
$ws.Range("D32").Value = "-"
$ws.Range("E32").Value = "-"
$ws.Range("F32").Value = "-"
$ws.Range("G32").Value = "-"
